$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# Revert the "survey" sheet's third column header back from "label" to "message"
# (undoing the earlier "wrong-xlsform-col" merge).
$survey.Range("C1").Value = "message"

# Restore the previously active/selected cell on the "survey" sheet.
$survey.Range("A3").Select()

# Cosmetic page-setup tweak that was part of the original revert.
$choices.PageSetup.FirstPageNumber = 1
